$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2")
$ws.Activate()

$rng = $ws.Range("F2:F34")
$rng.Value = $true
$rng.NumberFormat = """TRUE"";""TRUE"";""FALSE"""
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108
$rng.WrapText = $true
